$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 199 (old rows 199-200 shift down to 201-202)
$ws.Rows.Item(199).Resize(2).Insert()

# New row 199: updated "Camote" entry (date/quality/prices/origin changed)
$ws.Cells.Item(199, 1).Value = 5
$ws.Cells.Item(199, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(199, 3).Value = "Maule"
$ws.Cells.Item(199, 4).Value = 44628
$ws.Cells.Item(199, 5).Value = 7
$ws.Cells.Item(199, 6).Value = 100112045
$ws.Cells.Item(199, 7).Value = "Zapallo"
$ws.Cells.Item(199, 8).Value = "Camote"
$ws.Cells.Item(199, 9).Value = "1a (cosecha)"
$ws.Cells.Item(199, 10).Value = 900
$ws.Cells.Item(199, 11).Value = 300
$ws.Cells.Item(199, 12).Value = 300
$ws.Cells.Item(199, 13).Value = 300
$ws.Cells.Item(199, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(199, 15).Value = "Región del Maule"
$ws.Cells.Item(199, 16).Value = 300
$ws.Cells.Item(199, 17).Value = 1
$ws.Cells.Item(199, 18).Value = "Hortaliza"

# New row 200: new "Paine" entry
$ws.Cells.Item(200, 1).Value = 5
$ws.Cells.Item(200, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(200, 3).Value = "Maule"
$ws.Cells.Item(200, 4).Value = 44628
$ws.Cells.Item(200, 5).Value = 7
$ws.Cells.Item(200, 6).Value = 100112045
$ws.Cells.Item(200, 7).Value = "Zapallo"
$ws.Cells.Item(200, 8).Value = "Paine"
$ws.Cells.Item(200, 9).Value = "1a (cosecha)"
$ws.Cells.Item(200, 10).Value = 2000
$ws.Cells.Item(200, 11).Value = 120
$ws.Cells.Item(200, 12).Value = 120
$ws.Cells.Item(200, 13).Value = 120
$ws.Cells.Item(200, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(200, 15).Value = "Región del Maule"
$ws.Cells.Item(200, 16).Value = 120
$ws.Cells.Item(200, 17).Value = 1
$ws.Cells.Item(200, 18).Value = "Hortaliza"
